# Auto-generated Excel COM-interop script
# Applies the cryptos.xlsx price/volume/name refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed by Excel as a number
# (e.g. "0.456", "1.00") are written as literal text by temporarily forcing
# a Text number format, then ClearFormats() restores the original (default) style
# so only the cell VALUE changes -- matching the source diff which touches no styles.

$ws.Range("D2").Value = '67.708.23'
$ws.Range("E2").Value = '  -0.38%  '
$ws.Range("D3").Value = '3.771.19'
$ws.Range("E3").Value = '  -1.71%  '
$ws.Range("E4").Value = '  +0.30%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.67'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '169.20'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +1.16%  '
$ws.Range("D7").Value = '3.759.22'
$ws.Range("E7").Value = '  -2.16%  '
$ws.Range("E8").Value = '  +0.14%  '
$ws.Range("E9").Value = '  +0.26%  '
$ws.Range("E10").Value = '  +1.56%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.50'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.90%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.456'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.24%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000273'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +5.69%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.85'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.15%  '
$ws.Range("D15").Value = '4.415.99'
$ws.Range("E15").Value = '  -1.41%  '
$ws.Range("D16").Value = '3.814.35'
$ws.Range("E16").Value = '  -0.31%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.91'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +5.49%  '
$ws.Range("D18").Value = '67.798.83'
$ws.Range("E18").Value = '  -0.41%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.26'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.54%  '
$ws.Range("E20").Value = '  +0.89%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.54'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.06%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '468.43'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.89%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.724'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.67%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000150'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -4.79%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.81'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.22%  '
$ws.Range("E26").Value = '  +1.38%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.17'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.52%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.34'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +4.05%  '
$ws.Range("E29").Value = '  +0.13%  '
$ws.Range("E30").Value = '  -1.43%  '
$ws.Range("D31").Value = '3.932.64'
$ws.Range("E31").Value = '  -1.32%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.25'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.70%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '30.42'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -2.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.17'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -2.20%  '
$ws.Range("D36").Value = '3.745.86'
$ws.Range("E36").Value = '  -1.49%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.82'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +6.74%  '
$ws.Range("E38").Value = '  +1.33%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.90'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.49%  '
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.138'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -1.82%  '
$ws.Range("B41").Value = 'Mantle'
$ws.Range("C41").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -2.00%  '
$ws.Range("E42").Value = '  +0.03%  '
$ws.Range("E43").Value = '  +1.57%  '
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("B45").Value = 'Stacks'
$ws.Range("C45").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.97'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.36%  '
$ws.Range("B46").Value = 'Cosmos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.73'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +2.00%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '46.30'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.50%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '403.84'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -4.34%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.000277'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -5.44%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '141.77'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.69%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0356'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.20%  '
